# --- Journal.xlsx refactor: extend "Parameters" sheet with WSD/model/context-window
# --- columns, fill in the "test" row's clustering params, and add a new
# --- "Time + Size" sheet for the upcoming timing/size experiment pipeline.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New columns P:Q on "Parameters" -------------------------------------

# Header group "WSD" spanning P4:Q4 (mirrors the B4:F4 / J4:M4 grouped headers)
$ws1.Range("B4").Copy()
$ws1.Range("P4").PasteSpecial(-4122)
$ws1.Range("F4").Copy()
$ws1.Range("Q4").PasteSpecial(-4122)
$ws1.Range("P4").Value = "WSD"

# Sub-headers
$ws1.Range("O5").Copy()
$ws1.Range("P5:Q5").PasteSpecial(-4122)
$ws1.Range("P5").Value = "model"
$ws1.Range("Q5").Value = "context window"

# "1MB" row data
$ws1.Range("P6").Value = "argmax_i(prod_k(cossim(s_i,c_k))))"
$ws1.Range("Q6").Value = 10

# "test" row: fill in the clustering algorithm columns (same as "1MB" row)
$ws1.Range("J7").Value = "cw"
$ws1.Range("K7").Value = 200
$ws1.Range("L7").Value = 200
$ws1.Range("M7").Value = 5

# New column widths for P and Q
$ws1.Columns.Item(16).ColumnWidth = 13.416666666666666
$ws1.Columns.Item(17).ColumnWidth = 14.166666666666666

# Restore the working selection on Parameters
$ws1.Range("M15").Select()

# --- New "Time + Size" sheet ----------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Time + Size"

# Match the workbook's usual (inch-based) page margins instead of the
# add-sheet default (0.7/0.75in) so the new sheet is consistent with Parameters
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

$ws2.Range("F28").Select()
